$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "2025-03-28"
$ws.Range("A35").ClearFormats()
$ws.Range("B35").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C35").Value = "NA"
$ws.Range("D35").Value = 1
